$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Standorte")

$ws.Range("C2:C6").Value = 0.1
